$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.317.94'
$ws.Range("E2").Value = '  +1.47%  '
$ws.Range("D3").Value = '3.596.65'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''244.18'
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("D6").Value = '''1.79'
$ws.Range("E6").Value = '  +16.98%  '
$ws.Range("D7").Value = '''652.94'
$ws.Range("D8").Value = '''0.425'
$ws.Range("E8").Value = '  +5.22%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +2.19%  '
$ws.Range("D11").Value = '3.593.84'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '''44.83'
$ws.Range("E12").Value = '  +3.99%  '
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").Value = '''6.47'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '4.266.65'
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").Value = '97.186.82'
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '3.588.62'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").Value = '''18.34'
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("D22").Value = '''0.529'
$ws.Range("E22").Value = '  +6.10%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '''518.12'
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '''3.48'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("E25").Value = '  +4.05%  '
$ws.Range("D26").Value = '''6.97'
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("D27").Value = '''103.27'
$ws.Range("E27").Value = '  +7.42%  '
$ws.Range("D28").Value = '''13.29'
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").Value = '''0.184'
$ws.Range("E29").Value = '  +25.30%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''2.99'
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''12.05'
$ws.Range("E31").Value = '  +3.73%  '
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").Value = '''0.190'
$ws.Range("E33").Value = '  +6.09%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '''31.87'
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = '''1.71'
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").Value = '''0.584'
$ws.Range("E37").Value = '  +3.15%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = '''8.80'
$ws.Range("E38").Value = '  +1.56%  '
$ws.Range("D39").Value = '''616.73'
$ws.Range("E39").Value = '  +2.81%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.154'
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").Value = '''1.93'
$ws.Range("E41").Value = '  +2.23%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.460'
$ws.Range("E42").Value = '  +40.65%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '''0.931'
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '''6.10'
$ws.Range("E45").Value = '  +5.08%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0450'
$ws.Range("E46").Value = '  +7.09%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '''2.34'
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").Value = '''23.66'
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '''8.69'
$ws.Range("E49").Value = '  +5.04%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''33.00'
$ws.Range("E50").Value = '  -4.35%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '''3.30'
$ws.Range("E51").Value = '  +6.50%  '
